$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.5114799999999484
$ws.Range("D2").Value = 0.3034800000000131
$ws.Range("E2").Value = 0.1105599999999962
$ws.Range("F2").Value = 0.05015999999999863
$ws.Range("G2").Value = 0.01667999999999999
$ws.Range("H2").Value = 0.007640000000000017
# Row 3
$ws.Range("C3").Value = 0.05527999999999842
$ws.Range("D3").Value = 0.3365600000000003
$ws.Range("E3").Value = 0.2623600000000291
$ws.Range("F3").Value = 0.1560800000000051
$ws.Range("G3").Value = 0.09387999999999685
$ws.Range("H3").Value = 0.09583999999999677
# Row 4
$ws.Range("C4").Value = 0.3853999999999814
$ws.Range("D4").Value = 0.3591199999999916
$ws.Range("E4").Value = 0.1462000000000021
$ws.Range("F4").Value = 0.06935999999999785
$ws.Range("G4").Value = 0.02579999999999962
$ws.Range("H4").Value = 0.01412000000000003
# Row 5
$ws.Range("C5").Value = 0.3273200000000039
$ws.Range("D5").Value = 0.352839999999994
$ws.Range("E5").Value = 0.1679200000000087
$ws.Range("F5").Value = 0.08823999999999708
$ws.Range("G5").Value = 0.03891999999999909
$ws.Range("H5").Value = 0.02475999999999966
# Row 7
$ws.Range("C7").Value = 0.08203999999999734
$ws.Range("D7").Value = 0.4348399999999622
$ws.Range("E7").Value = 0.2377200000000301
$ws.Range("F7").Value = 0.1390399999999999
$ws.Range("G7").Value = 0.06199999999999815
$ws.Range("H7").Value = 0.04435999999999887
# Row 8
$ws.Range("C8").Value = 0.7786400000002156
$ws.Range("D8").Value = 0.1595600000000062
$ws.Range("E8").Value = 0.04279999999999893
$ws.Range("F8").Value = 0.01424000000000003
$ws.Range("G8").Value = 0.003400000000000006
$ws.Range("H8").Value = 0.001360000000000001
# Row 9
$ws.Range("C9").Value = 0.05119999999999859
$ws.Range("D9").Value = 0.3251600000000047
$ws.Range("E9").Value = 0.2188800000000243
$ws.Range("F9").Value = 0.1645200000000077
$ws.Range("G9").Value = 0.1010799999999966
$ws.Range("H9").Value = 0.1391599999999999
# Row 10
$ws.Range("C10").Value = 0.4495199999999565
$ws.Range("D10").Value = 0.3697599999999874
$ws.Range("E10").Value = 0.1168799999999959
$ws.Range("F10").Value = 0.04543999999999882
$ws.Range("G10").Value = 0.01284000000000003
$ws.Range("H10").Value = 0.005560000000000012
# Row 12
$ws.Range("C12").Value = 0.02651999999999959
$ws.Range("D12").Value = 0.4758799999999463
$ws.Range("E12").Value = 0.2496800000000338
$ws.Range("F12").Value = 0.1430000000000011
$ws.Range("G12").Value = 0.06083999999999819
$ws.Range("H12").Value = 0.04407999999999888
# Row 13
$ws.Range("C13").Value = 0.2194800000000245
$ws.Range("D13").Value = 0.4835199999999433
$ws.Range("E13").Value = 0.1836800000000136
$ws.Range("F13").Value = 0.07383999999999767
$ws.Range("G13").Value = 0.02607999999999961
$ws.Range("H13").Value = 0.01340000000000003
# Row 14
$ws.Range("C14").Value = 0.4346799999999623
$ws.Range("D14").Value = 0.3903599999999794
$ws.Range("E14").Value = 0.1188799999999958
$ws.Range("F14").Value = 0.03959999999999906
$ws.Range("G14").Value = 0.01164000000000003
$ws.Range("H14").Value = 0.00484000000000001
# Row 15
$ws.Range("C15").Value = 0.7468800000001838
$ws.Range("D15").Value = 0.1993200000000183
$ws.Range("E15").Value = 0.04123999999999899
$ws.Range("F15").Value = 0.01040000000000002
$ws.Range("G15").Value = 0.001680000000000002
$ws.Range("H15").Value = 0.00048
# Row 17
$ws.Range("C17").Value = 0.02207999999999977
$ws.Range("D17").Value = 0.4501999999999562
$ws.Range("E17").Value = 0.2456400000000325
$ws.Range("F17").Value = 0.1066799999999963
$ws.Range("G17").Value = 0.0805599999999974
$ws.Range("H17").Value = 0.09483999999999682
# Row 18
$ws.Range("C18").Value = 0.3681199999999881
$ws.Range("D18").Value = 0.3618399999999905
$ws.Range("E18").Value = 0.1703600000000095
$ws.Range("F18").Value = 0.06371999999999808
$ws.Range("G18").Value = 0.0238399999999997
$ws.Range("H18").Value = 0.01212000000000003
# Row 19
$ws.Range("C19").Value = 0.4899999999999408
$ws.Range("D19").Value = 0.3513599999999946
$ws.Range("E19").Value = 0.1199199999999958
$ws.Range("F19").Value = 0.02811999999999953
$ws.Range("G19").Value = 0.008000000000000018
$ws.Range("H19").Value = 0.002600000000000004
# Row 20
$ws.Range("C20").Value = 0.4907999999999405
$ws.Range("D20").Value = 0.3157200000000084
$ws.Range("E20").Value = 0.1318799999999977
$ws.Range("F20").Value = 0.04307999999999892
$ws.Range("G20").Value = 0.01280000000000003
$ws.Range("H20").Value = 0.005720000000000012
# Row 22
$ws.Range("C22").Value = 0.01416000000000003
$ws.Range("D22").Value = 0.298560000000015
$ws.Range("E22").Value = 0.3163600000000081
$ws.Range("F22").Value = 0.1530400000000042
$ws.Range("G22").Value = 0.0976399999999967
$ws.Range("H22").Value = 0.1202399999999958
# Row 23
$ws.Range("C23").Value = 0.234160000000029
$ws.Range("D23").Value = 0.3743199999999857
$ws.Range("E23").Value = 0.2421600000000315
$ws.Range("F23").Value = 0.08851999999999707
$ws.Range("G23").Value = 0.03919999999999908
$ws.Range("H23").Value = 0.02163999999999979
# Row 24
$ws.Range("C24").Value = 0.4029599999999746
$ws.Range("D24").Value = 0.3311200000000024
$ws.Range("E24").Value = 0.1702400000000094
$ws.Range("F24").Value = 0.06219999999999814
$ws.Range("G24").Value = 0.02207999999999977
$ws.Range("H24").Value = 0.01140000000000003
# Row 25
$ws.Range("C25").Value = 0.74300000000018
$ws.Range("D25").Value = 0.1714800000000098
$ws.Range("E25").Value = 0.06391999999999808
$ws.Range("F25").Value = 0.01652
$ws.Range("G25").Value = 0.004000000000000008
$ws.Range("H25").Value = 0.00108
# Row 27
$ws.Range("C27").Value = 0.7399600000001769
$ws.Range("D27").Value = 0.1672800000000085
$ws.Range("E27").Value = 0.06047999999999821
$ws.Range("F27").Value = 0.02311999999999973
$ws.Range("G27").Value = 0.006520000000000014
$ws.Range("H27").Value = 0.002640000000000004
# Row 28
$ws.Range("C28").Value = 0.1280399999999965
$ws.Range("D28").Value = 0.3224800000000058
$ws.Range("E28").Value = 0.2056400000000203
$ws.Range("F28").Value = 0.149120000000003
$ws.Range("G28").Value = 0.090239999999997
$ws.Range("H28").Value = 0.1044799999999964
# Row 29
$ws.Range("C29").Value = 0.09147999999999695
$ws.Range("D29").Value = 0.2789200000000226
$ws.Range("E29").Value = 0.1894800000000153
$ws.Range("F29").Value = 0.1520400000000039
$ws.Range("G29").Value = 0.1038399999999964
$ws.Range("H29").Value = 0.1842400000000137
# Row 30
$ws.Range("C30").Value = 0.2068000000000206
$ws.Range("D30").Value = 0.386319999999981
$ws.Range("E30").Value = 0.1871200000000146
$ws.Range("F30").Value = 0.1199199999999958
$ws.Range("G30").Value = 0.05643999999999837
$ws.Range("H30").Value = 0.0433999999999989

Write-Output "Updated 144 probability cells (C2:H30) with new simulation results."
